# Tidsplan_xjobb.xlsx - "Solved main post to db issue for now"
#
# The react-app / db update task (row 19) finally went through, so:
#   - ActlHours (P19)        25   -> 32   (more hours were actually spent)
#   - ActDeliveryDate (Q19)  "2020-02-27 Pending" -> "2020-02-28 Pending"
#   - Comment (T19)          "Greate trouble making it work to update db
#                              via react app"
#                             -> "Greate trouble making it work to update db
#                              via react app. Finally got through.."
#
# All of the Plan-vs-Act / summary formulas further down the sheet (S19,
# P36, S36, P37, S37, P38, P43, P44, P45) are plain formulas that depend on
# these inputs, so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P19").Value = 32
$ws.Range("Q19").Value = "2020-02-28 Pending"
$ws.Range("T19").Value = "Greate trouble making it work to update db via react app. Finally got through.."

# Leave the cursor where the author left it when saving.
[void]$ws.Range("T19").Select()
